# Update the "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Updates common to both sheets (展览 and 全部类型)
$updates = @{
    3  = 92
    5  = 48
    6  = 559
    8  = 2023
    11 = 4360
    14 = 100
    15 = 7
    16 = 113
    17 = 25
    19 = 70
    20 = 3157
    21 = 68
    22 = 465
    25 = 74
    29 = 54
    32 = 534
    34 = 268
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}
$ws1.Range("F33").Value = 1736

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $ws4.Range("F$row").Value = $updates[$row]
}
$ws4.Range("F33").Value = 1737
